# Update the cached "datetimeFigureOut" date field text wherever it
# appears (slide master + every slide layout) from "10/01/2023" to
# "2023-10-05".
$p = $ppt.ActivePresentation

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame) {
                if ($sh.TextFrame.HasText) {
                    $sh.TextFrame.TextRange.Text = "2023-10-05"
                }
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster
for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    Update-DatePlaceholder $p.SlideMaster.CustomLayouts.Item($j)
}

# Slide content: swap the "c" / "d" labels between the two ovals and
# nudge "Oval 34" up slightly.
$s = $p.Slides.Item(1)

$ovalC = $s.Shapes.Item("Oval 32")   # currently labelled "c"
$ovalD = $s.Shapes.Item("Oval 34")   # currently labelled "d"

$ovalC.TextFrame.TextRange.Text = "d"
$ovalD.TextFrame.TextRange.Text = "c"
$ovalD.Top = 42.46669291338583
